$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing "Next_update" dates ---
$ws.Range("D4").Value = 45261
$ws.Range("D5").Value = 45261

# Row 6's Next_update (D6) switches from the datetime-with-time style to the
# plain yyyy-mm-dd date style (matching the style already used by column D
# elsewhere), then gets the same updated date.
$ws.Range("D6").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("D6").Value = 45261

# --- Append new vendor row: TP-Link ---
$ws.Range("A7").Value = "TP-Link"
$ws.Range("B7").Value = 0
$ws.Range("C7:D7").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("C7").Value = 44902
$ws.Range("D7").Value = 44902
$ws.Range("E7").Value = "TPLinkScraper"
$ws.Range("F7").Value = 20

# --- Column E gets an explicit width to fit the new scraper-class values ---
$ws.Columns.Item(5).ColumnWidth = 29.166666666666668

# --- Selection moves as the user continued working below the table ---
$ws.Range("E19").Select()
